$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$shp = $s.Shapes.Item(1)
$tbl = $shp.Table

# Distribute the heights of the last three rows (rows 4-6) evenly,
# mirroring PowerPoint's "Distribute Rows" command. The new, equal
# row height (in points) is derived from the total EMU height of the
# three source rows divided evenly across them.
$newHeightPoints = 379723 / 12700

$tbl.Rows.Item(4).Height = $newHeightPoints
$tbl.Rows.Item(5).Height = $newHeightPoints
$tbl.Rows.Item(6).Height = $newHeightPoints
